# Auto-generated Excel COM-interop script
# Applies the "Phantom_Profits" value-refresh diff (scheduled runner update)
# across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 9722
$ws.Range("I34").Value = 6296
$ws.Range("K34").Value = 6296
$ws.Range("M34").Value = -6093
$ws.Range("H36").Value = 9722
$ws.Range("I36").Value = 6296
$ws.Range("K36").Value = 6296
$ws.Range("M36").Value = -5581
$ws.Range("H40").Value = 3790
$ws.Range("I40").Value = 1119.8
$ws.Range("J40").Value = 4902.5835
$ws.Range("K40").Value = 1119.8
$ws.Range("L40").Value = 4902.5835
$ws.Range("M40").Value = -944.8
$ws.Range("N40").Value = -5252.5835
$ws.Range("H53").Value = 483.75
$ws.Range("I53").Value = 437.77777
$ws.Range("K53").Value = 437.77777
$ws.Range("M53").Value = 199.22223
$ws.Range("H58").Value = 483.8
$ws.Range("I58").Value = 140
$ws.Range("J58").Value = 999.5
$ws.Range("K58").Value = 420
$ws.Range("L58").Value = 2998.5
$ws.Range("M58").Value = -270
$ws.Range("N58").Value = -3298.5
$ws.Range("H92").Value = 294.83334
$ws.Range("I92").Value = 294.83334
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 294.83334
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 953.16666
$ws.Range("H135").Value = 469.33334
$ws.Range("I135").Value = 464.69232
$ws.Range("K135").Value = 4182.23088
$ws.Range("M135").Value = -1647.23088
$ws.Range("H137").Value = 2861.5652
$ws.Range("J137").Value = 3165.7778
$ws.Range("L137").Value = 9497.3334
$ws.Range("N137").Value = -14597.3334
$ws.Range("H138").Value = 3080.861
$ws.Range("J138").Value = 3998.5454
$ws.Range("L138").Value = 11995.6362
$ws.Range("N138").Value = -22275.6362
$ws.Range("H140").Value = 90780
$ws.Range("J140").Value = 90780
$ws.Range("L140").Value = 90780
$ws.Range("N140").Value = -101140
$ws.Range("N92").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5166.6665
$ws.Range("I39").Value = 5166.6665
$ws.Range("K39").Value = 5166.6665
$ws.Range("M39").Value = -4646.6665
$ws.Range("H61").Value = 1875.2
$ws.Range("I61").Value = 1861.3334
$ws.Range("K61").Value = 1861.3334
$ws.Range("M61").Value = -1649.3334
$ws.Range("H74").Value = 2506.7083
$ws.Range("I74").Value = 2591.0476
$ws.Range("J74").Value = 1916.3334
$ws.Range("K74").Value = 2591.0476
$ws.Range("L74").Value = 1916.3334
$ws.Range("M74").Value = -1717.0476
$ws.Range("N74").Value = -3664.3334
$ws.Range("H77").Value = 2506.7083
$ws.Range("I77").Value = 2591.0476
$ws.Range("J77").Value = 1916.3334
$ws.Range("K77").Value = 12955.238
$ws.Range("L77").Value = 9581.666999999999
$ws.Range("M77").Value = -8587.237999999999
$ws.Range("N77").Value = -18317.667
$ws.Range("H110").Value = 3130.1765
$ws.Range("I110").Value = 3130.1765
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3130.1765
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1085.1765
$ws.Range("H132").Value = 3690.8
$ws.Range("I132").Value = 3690.8
$ws.Range("K132").Value = 11072.4
$ws.Range("M132").Value = -8542.400000000001
$ws.Range("H136").Value = 1875.2
$ws.Range("I136").Value = 1861.3334
$ws.Range("K136").Value = 5584.0002
$ws.Range("M136").Value = -3034.0002
$ws.Range("N110").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("H43").Value = 9995
$ws.Range("J43").Value = 9995
$ws.Range("L43").Value = 9995
$ws.Range("N43").Value = -10363
$ws.Range("H47").Value = 8429.666999999999
$ws.Range("I47").Value = 8429.666999999999
$ws.Range("K47").Value = 8429.666999999999
$ws.Range("M47").Value = -7863.666999999999
$ws.Range("H55").Value = 8899.666999999999
$ws.Range("I55").Value = 8899.666999999999
$ws.Range("K55").Value = 8899.666999999999
$ws.Range("M55").Value = -8584.666999999999
$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101372
$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306864
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("H101").Value = 9995
$ws.Range("J101").Value = 9995
$ws.Range("L101").Value = 9995
$ws.Range("N101").Value = -16485
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1271814.2
$ws.Range("I4").Value = 660576.3
$ws.Range("K4").Value = 1981728.9
$ws.Range("M4").Value = -1981616.9
$ws.Range("H11").Value = 972.5714
$ws.Range("I11").Value = 1109.3334
$ws.Range("J11").Value = 152
$ws.Range("K11").Value = 3328.0002
$ws.Range("L11").Value = 456
$ws.Range("M11").Value = -3188.0002
$ws.Range("N11").Value = -736
$ws.Range("H37").Value = 250000
$ws.Range("J37").Value = 250000
$ws.Range("L37").Value = 750000
$ws.Range("N37").Value = -750224
$ws.Range("H41").Value = 200.5
$ws.Range("I41").Value = 200.5
$ws.Range("K41").Value = 601.5
$ws.Range("M41").Value = -263.5
$ws.Range("H116").Value = 12032
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 12032
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 36096
$ws.Range("N116").Value = -42980
$ws.Range("M116").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7998
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 7998
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164
$ws.Range("H101").Value = 13537.6
$ws.Range("J101").Value = 13537.6
$ws.Range("L101").Value = 13537.6
$ws.Range("N101").Value = -20027.6
$ws.Range("H126").Value = 2174
$ws.Range("I126").Value = 2174
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6522
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4052
$ws.Range("H132").Value = 1627.4445
$ws.Range("I132").Value = 1627.4445
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4882.333500000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2352.333500000001
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("N132").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 520.375
$ws.Range("I22").Value = 494
$ws.Range("J22").Value = 599.5
$ws.Range("K22").Value = 494
$ws.Range("L22").Value = 599.5
$ws.Range("M22").Value = -199
$ws.Range("N22").Value = -1189.5
$ws.Range("H27").Value = 520.375
$ws.Range("I27").Value = 494
$ws.Range("J27").Value = 599.5
$ws.Range("K27").Value = 494
$ws.Range("L27").Value = 599.5
$ws.Range("M27").Value = -387
$ws.Range("N27").Value = -813.5
$ws.Range("H54").Value = 14350
$ws.Range("J54").Value = 14350
$ws.Range("L54").Value = 14350
$ws.Range("N54").Value = -15638
$ws.Range("H132").Value = 1692.3462
$ws.Range("I132").Value = 1170.9166
$ws.Range("K132").Value = 3512.7498
$ws.Range("M132").Value = -982.7498000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3773.7307
$ws.Range("I136").Value = 4105.905
$ws.Range("J136").Value = 2378.6
$ws.Range("K136").Value = 12317.715
$ws.Range("L136").Value = 7135.799999999999
$ws.Range("M136").Value = -9767.715
$ws.Range("N136").Value = -12235.8
